$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 455, pushing existing rows 455:503 down to 456:504
$ws.Rows.Item(455).Insert()

# Populate the newly inserted row 455 with the new record
$ws.Cells.Item(455, 1).Value = 10
$ws.Cells.Item(455, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(455, 3).Value = "La Araucanía"
$ws.Cells.Item(455, 4).Value = 44946
$ws.Cells.Item(455, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(455, 5).Value = 9
$ws.Cells.Item(455, 6).Value = 100112040
$ws.Cells.Item(455, 7).Value = "Cilantro"
$ws.Cells.Item(455, 8).Value = "Sin especificar"
$ws.Cells.Item(455, 9).Value = "Primera"
$ws.Cells.Item(455, 10).Value = 95
$ws.Cells.Item(455, 11).Value = 5000
$ws.Cells.Item(455, 12).Value = 6000
$ws.Cells.Item(455, 13).Value = 5368
$ws.Cells.Item(455, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(455, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(455, 16).Value = 2684
$ws.Cells.Item(455, 17).Value = 2
$ws.Cells.Item(455, 18).Value = "Hortaliza"
